$d = $word.ActiveDocument

# Locate the paragraph that holds "Docente(s) Responsável(eis) " by index,
# so we can insert a brand new paragraph right after it.
$targetIndex = 0
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text -like "*Docente(s) Responsável(eis)*") {
        $targetIndex = $i
    }
}

if ($targetIndex -gt 0) {
    # Insert a new (empty) paragraph right after the target paragraph.
    $d.Paragraphs($targetIndex).Range.InsertParagraphAfter()

    # The newly created paragraph is now the next one; give it the
    # ListBullet style and the docente's name/id text.
    $newPara = $d.Paragraphs($targetIndex + 1)
    $newPara.Style = "ListBullet"
    $newPara.Range.Text = "7455355 - Robson da Silva Rocha"
}
